$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "320.24"
Set-TextValue "E2" "4.92%"
Set-TextValue "G2" "20"

Set-TextValue "D3" "49.70"
Set-TextValue "E3" "12.33%"
Set-TextValue "G3" "20"

Set-TextValue "D4" "5.335"
Set-TextValue "E4" "4.18%"
Set-TextValue "G4" "20"

Set-TextValue "D5" "0.08022"
Set-TextValue "E5" "2.25%"
Set-TextValue "G5" "20"

Set-TextValue "D6" "4.598"
Set-TextValue "E6" "3.61%"
Set-TextValue "G6" "20"

Set-TextValue "D7" "1.341"
Set-TextValue "E7" "26.89%"
Set-TextValue "G7" "20"

Set-TextValue "D8" "1.647"
Set-TextValue "E8" "1.55%"
Set-TextValue "G8" "20"

Set-TextValue "D9" "0.1253"
Set-TextValue "E9" "-4.52%"
Set-TextValue "G9" "20"

Set-TextValue "D10" "0.1964"
Set-TextValue "E10" "4.98%"
Set-TextValue "G10" "20"

Set-TextValue "D11" "0.09669"
Set-TextValue "E11" "5.30%"
Set-TextValue "G11" "20"

Set-TextValue "D12" "0.04536"
Set-TextValue "E12" "8.72%"
Set-TextValue "G12" "20"

Set-TextValue "D13" "0.1047"
Set-TextValue "E13" "0.25%"
Set-TextValue "G13" "20"

Set-TextValue "D14" "0.001315"
Set-TextValue "E14" "2.06%"
Set-TextValue "G14" "20"

Set-TextValue "E15" "0.99%"
Set-TextValue "G15" "20"

Set-TextValue "D16" "0.005855"
Set-TextValue "E16" "2.03%"
Set-TextValue "G16" "20"

Set-TextValue "E17" "-0.67%"
Set-TextValue "G17" "20"

Set-TextValue "E18" "5.71%"
Set-TextValue "G18" "20"

Set-TextValue "D19" "0.3477"
Set-TextValue "E19" "3.07%"
Set-TextValue "G19" "20"

Set-TextValue "D20" "8.179"
Set-TextValue "E20" "1.59%"
Set-TextValue "G20" "20"

Set-TextValue "D21" "0.1392"
Set-TextValue "E21" "0.74%"
Set-TextValue "G21" "20"

Set-TextValue "D22" "0.3012"
Set-TextValue "E22" "7.26%"
Set-TextValue "G22" "20"

Set-TextValue "E23" "1.46%"
Set-TextValue "G23" "20"

Set-TextValue "D24" "0.004225"
Set-TextValue "E24" "-5.90%"
Set-TextValue "G24" "20"

Set-TextValue "E25" "0.67%"
Set-TextValue "G25" "20"

Set-TextValue "D26" "0.0003543"
Set-TextValue "G26" "20"

Set-TextValue "G27" "20"

Set-TextValue "G28" "20"

Set-TextValue "G29" "20"

Set-TextValue "G30" "20"

Set-TextValue "G31" "20"

Set-TextValue "G32" "20"

Set-TextValue "G33" "20"

Set-TextValue "G34" "20"

Set-TextValue "G35" "20"

Set-TextValue "G36" "20"

Set-TextValue "G37" "20"

Set-TextValue "D38" "0.02666"
Set-TextValue "E38" "1.74%"
Set-TextValue "G38" "20"

Set-TextValue "D39" "0.05895"
Set-TextValue "E39" "10.18%"
Set-TextValue "G39" "20"

Set-TextValue "E40" "92.88%"
Set-TextValue "G40" "20"

Set-TextValue "D41" "0.008045"
Set-TextValue "E41" "4.00%"
Set-TextValue "G41" "20"

Set-TextValue "D42" "0.1467"
Set-TextValue "E42" "6.29%"
Set-TextValue "G42" "20"

Set-TextValue "D43" "0.007524"
Set-TextValue "E43" "2.71%"
Set-TextValue "G43" "20"

Set-TextValue "D44" "0.007930"
Set-TextValue "E44" "-4.97%"
Set-TextValue "G44" "20"

Set-TextValue "D45" "0.3219"
Set-TextValue "E45" "6.53%"
Set-TextValue "G45" "20"

Set-TextValue "D46" "0.00007013"
Set-TextValue "G46" "20"

Set-TextValue "E47" "0.66%"
Set-TextValue "G47" "20"

Set-TextValue "D48" "0.05594"
Set-TextValue "E48" "15.97%"
Set-TextValue "G48" "20"

Set-TextValue "D49" "0.004004"
Set-TextValue "E49" "0.67%"
Set-TextValue "G49" "20"

Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.66%"
Set-TextValue "G50" "20"

Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.66%"
Set-TextValue "G51" "20"
